$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new leading data row for 2022-Q3,
#    pushing the existing quarters down by one row.
# -----------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Copy the format of the last existing data row down into the new last
# row (row 9) before we touch any values, so the new row 9 / column A
# cell picks up the same style (s="2") as the rest of the index column.
$summary.Cells.Item(8, 1).Copy()
$summary.Cells.Item(9, 1).PasteSpecial(-4122)

# Shift rows 2..8 down to 3..9 (copy from the bottom up so we never
# clobber a row before it has been read).
for ($r = 8; $r -ge 2; $r--) {
    $dst = $r + 1
    $summary.Cells.Item($dst, 2).Value = $summary.Cells.Item($r, 2).Value2
    $summary.Cells.Item($dst, 3).Value = $summary.Cells.Item($r, 3).Value2
    $summary.Cells.Item($dst, 4).Value = $summary.Cells.Item($r, 4).Value2
    $summary.Cells.Item($dst, 1).Value = $r - 1
}

# New first data row: 2022-Q3
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 5
$summary.Cells.Item(2, 4).Value = 0.21

# -----------------------------------------------------------------------
# 2) Add a brand-new worksheet "2022-Q3" right after "总计" (i.e. before
#    the sheet that used to be "2022-Q2"), holding the fund-holding
#    detail rows for the new quarter.
# -----------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q3"

# Match the page-margin convention used by the rest of the workbook
# (0.75in / 1in / 0.5in -> 54 / 72 / 36 points).
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Borrow the header / index-column formatting from the sheet that has
# the same layout (now shifted one slot to index 3 - "2022-Q2").
$refSheet = $wb.Worksheets.Item(3)
$refSheet.Range("B1:H1").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)
$refSheet.Range("A2:A6").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Header row
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# Index column (A2:A6)
$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(4, 1).Value = 2
$newSheet.Cells.Item(5, 1).Value = 3
$newSheet.Cells.Item(6, 1).Value = 4

# Force the fund-code / ratio columns to plain text so leading/trailing
# zeros in values such as "012971" or "91.10" survive unchanged.
$newSheet.Range("B2:G6").NumberFormat = "@"

$newSheet.Cells.Item(2, 2).Value = "377150"
$newSheet.Cells.Item(2, 3).Value = "上投摩根健康品质生活混合A"
$newSheet.Cells.Item(2, 4).Value = "2.45"
$newSheet.Cells.Item(2, 5).Value = "90.27"
$newSheet.Cells.Item(2, 6).Value = "3.16"
$newSheet.Cells.Item(2, 7).Value = "0.0774"
$newSheet.Cells.Item(2, 8).Value = 10

$newSheet.Cells.Item(3, 2).Value = "012971"
$newSheet.Cells.Item(3, 3).Value = "东吴消费成长混合A"
$newSheet.Cells.Item(3, 4).Value = "0.84"
$newSheet.Cells.Item(3, 5).Value = "91.10"
$newSheet.Cells.Item(3, 6).Value = "6.06"
$newSheet.Cells.Item(3, 7).Value = "0.0509"
$newSheet.Cells.Item(3, 8).Value = 5

$newSheet.Cells.Item(4, 2).Value = "015346"
$newSheet.Cells.Item(4, 3).Value = "上投摩根健康品质生活混合C"
$newSheet.Cells.Item(4, 4).Value = "1.50"
$newSheet.Cells.Item(4, 5).Value = "90.27"
$newSheet.Cells.Item(4, 6).Value = "3.16"
$newSheet.Cells.Item(4, 7).Value = "0.0474"
$newSheet.Cells.Item(4, 8).Value = 10

$newSheet.Cells.Item(5, 2).Value = "011389"
$newSheet.Cells.Item(5, 3).Value = "国都聚成混合"
$newSheet.Cells.Item(5, 4).Value = "0.42"
$newSheet.Cells.Item(5, 5).Value = "78.46"
$newSheet.Cells.Item(5, 6).Value = "3.94"
$newSheet.Cells.Item(5, 7).Value = "0.0165"
$newSheet.Cells.Item(5, 8).Value = 5

$newSheet.Cells.Item(6, 2).Value = "012972"
$newSheet.Cells.Item(6, 3).Value = "东吴消费成长混合C"
$newSheet.Cells.Item(6, 4).Value = "0.24"
$newSheet.Cells.Item(6, 5).Value = "91.10"
$newSheet.Cells.Item(6, 6).Value = "6.06"
$newSheet.Cells.Item(6, 7).Value = "0.0145"
$newSheet.Cells.Item(6, 8).Value = 5

Write-Host "2022-Q3 sheet inserted and 总计 updated"
